# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "52.298.20"
$ws.Range("E2").Value = "  +5.66%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.802.39"
$ws.Range("E3").Value = "  +6.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - Solana
$ws.Range("D5").Value = "117.21"
$ws.Range("E5").Value = "  +4.81%  "

# Row 6 - BNB
$ws.Range("D6").Value = "341.96"
$ws.Range("E6").Value = "  +4.91%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +5.39%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  +6.42%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "42.31"
$ws.Range("E10").Value = "  +7.11%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +7.56%  "

# Row 12 - Chainlink
$ws.Range("D12").Value = "20.18"
$ws.Range("E12").Value = "  +0.35%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.33%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.92%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.238.73"
$ws.Range("E15").Value = "  +6.17%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.819.42"
$ws.Range("E16").Value = "  +6.74%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  +4.57%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "52.205.61"
$ws.Range("E18").Value = "  +5.55%  "

# Row 19 - ImmutableX
$ws.Range("D19").Value = "3.23"
$ws.Range("E19").Value = "  +11.14%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "13.44"
$ws.Range("E20").Value = "  +1.56%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.97"
$ws.Range("E21").Value = "  +4.35%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0990"
$ws.Range("E22").Value = "  +4.51%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "280.12"
$ws.Range("E23").Value = "  +4.26%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "70.45"
$ws.Range("E24").Value = "  +2.02%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.83"
$ws.Range("E25").Value = "  +10.22%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "27.00"
$ws.Range("E26").Value = "  +3.61%  "

# Row 27 - Dai
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.60%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +1.17%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +3.06%  "

# Row 31 - InjectiveProtocol
$ws.Range("E31").Value = "  +1.36%  "

# Row 32 - OKB
$ws.Range("D32").Value = "50.45"
$ws.Range("E32").Value = "  +1.70%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +5.05%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.0829"
$ws.Range("E34").Value = "  +2.63%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "2.14"
$ws.Range("E35").Value = "  +5.24%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.02%  "

# Row 37 - Celestia
$ws.Range("D37").Value = "19.06"
$ws.Range("E37").Value = "  +0.27%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "5.02"
$ws.Range("E38").Value = "  +1.43%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  +5.76%  "

# Row 40 - Stacks
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  +28.39%  "

# Row 41 - VeChain
$ws.Range("D41").Value = "0.0371"
$ws.Range("E41").Value = "  +11.94%  "

# Row 42 - EnergySwap
$ws.Range("D42").Value = "23.60"
$ws.Range("E42").Value = "  +4.13%  "

# Row 43 & 44 swapped - Stellar now ranks above WEMIXToken
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  +4.16%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +4.91%  "

# Row 45 - Monero
$ws.Range("D45").Value = "126.44"
$ws.Range("E45").Value = "  -2.07%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.110.90"
$ws.Range("E46").Value = "  +2.33%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  +2.96%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  +3.42%  "

# Row 49 - THORChain
$ws.Range("E49").Value = "  +7.22%  "

# Row 50 - SEI
$ws.Range("D50").Value = "0.917"
$ws.Range("E50").Value = "  +22.45%  "

# Row 51 - FraxShare
$ws.Range("D51").Value = "8.98"
$ws.Range("E51").Value = "  +1.17%  "
